# Apply the "Статус" (Status) column updates + E2E description tweak
# to the "Тест-кейсы на автоматизацию" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Тест-кейсы на автоматизацию")

# Test Plans CRUD block (rows 2-5) -> "In Qase"
$ws.Range("D2").Value = "In Qase"
$ws.Range("D3").Value = "In Qase"
$ws.Range("D4").Value = "In Qase"
$ws.Range("D5").Value = "In Qase"

# Test Cases CRUD block (rows 7-10) -> "Automated"
$ws.Range("D7").Value = "Automated"
$ws.Range("D8").Value = "Automated"
$ws.Range("D9").Value = "Automated"
$ws.Range("D10").Value = "Automated"

# Login block (rows 12-13) -> "Automated"
$ws.Range("D12").Value = "Automated"
$ws.Range("D13").Value = "Automated"

# E2E description: "Create 3 Test Cases..." -> "Create 2 Test Cases..."
$ws.Range("B14").Value = "Create 2 Test Cases and add to Test Plan"

# Widen the Status column slightly to fit the new text
$ws.Columns("D").ColumnWidth = 10.8

# Update the active selection to match the author's final cursor position
$ws.Range("G11").Select()
